$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2150537634408602
$ws.Range("C2").Value = 0.5125448028673835
$ws.Range("J2").Value = 0.02867383512544803
$ws.Range("P2").Value = 0.1397849462365591
$ws.Range("S2").Value = 0.1039426523297491
$ws.Range("B3").Value = 0.006896551724137931
$ws.Range("C3").Value = 0.02758620689655172
$ws.Range("J3").Value = 0.02758620689655172
$ws.Range("P3").Value = 0.7517241379310344
$ws.Range("S3").Value = 0.1862068965517241
$ws.Range("J4").Value = 0.08771929824561403
$ws.Range("P4").Value = 0.5964912280701754
$ws.Range("S4").Value = 0.3157894736842105
$ws.Range("B6").Value = 0.04390243902439024
$ws.Range("D6").Value = 0.01951219512195122
$ws.Range("E6").Value = 0.004878048780487805
$ws.Range("F6").Value = 0.04390243902439024
$ws.Range("J6").Value = 0.2
$ws.Range("O6").Value = 0.02439024390243903
$ws.Range("Q6").Value = 0.2390243902439024
$ws.Range("R6").Value = 0.06829268292682927
$ws.Range("S6").Value = 0.3560975609756097
$ws.Range("B7").Value = 0.06986899563318777
$ws.Range("D7").Value = 0.02183406113537118
$ws.Range("E7").Value = 0.004366812227074236
$ws.Range("F7").Value = 0.06550218340611354
$ws.Range("J7").Value = 0.1004366812227074
$ws.Range("O7").Value = 0.008733624454148471
$ws.Range("Q7").Value = 0.2096069868995633
$ws.Range("R7").Value = 0.09606986899563319
$ws.Range("S7").Value = 0.4235807860262009
$ws.Range("B8").Value = 0.08829174664107485
$ws.Range("D8").Value = 0.01919385796545105
$ws.Range("F8").Value = 0.05374280230326296
$ws.Range("J8").Value = 0.08829174664107485
$ws.Range("O8").Value = 0.01919385796545105
$ws.Range("Q8").Value = 0.2111324376199616
$ws.Range("R8").Value = 0.07677543186180422
$ws.Range("S8").Value = 0.4433781190019194
$ws.Range("B9").Value = 0.1111111111111111
$ws.Range("D9").Value = 0.006172839506172839
$ws.Range("E9").Value = 0.006172839506172839
$ws.Range("F9").Value = 0.06790123456790123
$ws.Range("J9").Value = 0.09876543209876543
$ws.Range("O9").Value = 0.0308641975308642
$ws.Range("Q9").Value = 0.1728395061728395
$ws.Range("R9").Value = 0.06790123456790123
$ws.Range("S9").Value = 0.4382716049382716
$ws.Range("B10").Value = 0.103765690376569
$ws.Range("D10").Value = 0.03096234309623431
$ws.Range("E10").Value = 0.001673640167364017
$ws.Range("F10").Value = 0.07196652719665272
$ws.Range("J10").Value = 0.1112970711297071
$ws.Range("O10").Value = 0.01757322175732218
$ws.Range("Q10").Value = 0.2317991631799163
$ws.Range("R10").Value = 0.07280334728033473
$ws.Range("S10").Value = 0.3581589958158996
$ws.Range("F11").Value = 0.002890173410404624
$ws.Range("G11").Value = 0.1445086705202312
$ws.Range("J11").Value = 0.0838150289017341
$ws.Range("K11").Value = 0.1907514450867052
$ws.Range("L11").Value = 0.5606936416184971
$ws.Range("S11").Value = 0.01734104046242774
$ws.Range("G12").Value = 0.7563451776649747
$ws.Range("J12").Value = 0.2131979695431472
$ws.Range("K12").Value = 0.01015228426395939
$ws.Range("L12").Value = 0.01015228426395939
$ws.Range("S12").Value = 0.01015228426395939
$ws.Range("F13").Value = 0.01923076923076923
$ws.Range("G13").Value = 0.7307692307692307
$ws.Range("J13").Value = 0.2307692307692308
$ws.Range("S13").Value = 0.01923076923076923
$ws.Range("F15").Value = 0.009569377990430622
$ws.Range("H15").Value = 0.2009569377990431
$ws.Range("I15").Value = 0.06220095693779904
$ws.Range("J15").Value = 0.354066985645933
$ws.Range("K15").Value = 0.05263157894736842
$ws.Range("M15").Value = 0.01435406698564593
$ws.Range("O15").Value = 0.03349282296650718
$ws.Range("S15").Value = 0.2727272727272727
$ws.Range("F16").Value = 0.02793296089385475
$ws.Range("H16").Value = 0.1843575418994413
$ws.Range("I16").Value = 0.1005586592178771
$ws.Range("J16").Value = 0.4189944134078212
$ws.Range("K16").Value = 0.1284916201117318
$ws.Range("M16").Value = 0.01675977653631285
$ws.Range("O16").Value = 0.0446927374301676
$ws.Range("S16").Value = 0.0782122905027933
$ws.Range("F17").Value = 0.01574803149606299
$ws.Range("H17").Value = 0.2047244094488189
$ws.Range("I17").Value = 0.07874015748031496
$ws.Range("J17").Value = 0.4114173228346457
$ws.Range("K17").Value = 0.1102362204724409
$ws.Range("M17").Value = 0.02755905511811024
$ws.Range("O17").Value = 0.05905511811023622
$ws.Range("S17").Value = 0.09251968503937008
$ws.Range("F18").Value = 0.02923976608187134
$ws.Range("H18").Value = 0.2222222222222222
$ws.Range("I18").Value = 0.07602339181286549
$ws.Range("J18").Value = 0.3859649122807017
$ws.Range("K18").Value = 0.152046783625731
$ws.Range("M18").Value = 0.01754385964912281
$ws.Range("O18").Value = 0.02339181286549707
$ws.Range("S18").Value = 0.0935672514619883
$ws.Range("F19").Value = 0.008123476848090982
$ws.Range("H19").Value = 0.2502030869212023
$ws.Range("I19").Value = 0.06417546709991877
$ws.Range("J19").Value = 0.3517465475223396
$ws.Range("K19").Value = 0.1251015434606012
$ws.Range("M19").Value = 0.02437043054427295
$ws.Range("O19").Value = 0.07148659626320066
$ws.Range("S19").Value = 0.1047928513403737
